# models tipo y subtipo, backend editar_prod
# Row 5 of the "templates" sheet (editar_prod template) changes:
#   - column G (estado2) goes from "en proceso" to "listo"
#   - column H (comentarios) comment "G: no aparecen los ingredientes" is removed (cell cleared)
# Also the view scrolls back to the top (A1) and the active selection moves to G5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G5").Value = "listo"
$ws.Range("H5").ClearContents()

# Reset the window scroll position (removes topLeftCell="A4") and select G5,
# matching the saved view state in the target workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G5").Select()
